$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-7 and extend with new rows 8-21 ---
# Row 2
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = 'Tomar un punto $A$'
$ws.Cells.Item(2, 2).WrapText = $true
$ws.Cells.Item(2, 3).Value = 'Hipótesis'
$ws.Cells.Item(2, 4).Value = 'punto'
$ws.Cells.Item(2, 5).Value = '{''x'':-3,''y'':0,''nombre'':''A''}'

# Row 3
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = 'Tomar un punto $B$'
$ws.Cells.Item(3, 2).WrapText = $true
$ws.Cells.Item(3, 3).Value = 'Hipótesis'
$ws.Cells.Item(3, 4).Value = 'punto'
$ws.Cells.Item(3, 5).Value = '{''x'':3,''y'':0,''nombre'':''B''}'

# Row 4
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = 'Trazar segmento infinito $\overline{AB}$ '
$ws.Cells.Item(4, 2).WrapText = $true
$ws.Cells.Item(4, 3).Value = 'Hipótesis'
$ws.Cells.Item(4, 4).Value = 'segmento'
$ws.Cells.Item(4, 5).Value = '{''x1'':-10, ''y1'':0, ''x2'':10, ''y2'':0}'

# Row 5
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = 'Tomar un punto $C$ fuera de la recta $\overline{AB}$'
$ws.Cells.Item(5, 2).WrapText = $true
$ws.Cells.Item(5, 3).Value = 'Hipótesis'
$ws.Cells.Item(5, 4).Value = 'punto'
$ws.Cells.Item(5, 5).Value = '{''x'':0,''y'':1,''nombre'':''C''}'

# Row 6
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = 'Tomar punto $D$ en el lado opuesto de $\overline{AB}$'
$ws.Cells.Item(6, 2).WrapText = $true
$ws.Cells.Item(6, 3).Value = '-'
$ws.Cells.Item(6, 4).Value = 'punto'
$ws.Cells.Item(6, 5).Value = '{''x'':1.5,''y'':-1,''nombre'':''D''}'

# Row 7
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = 'Trazar circunferencia con centro $C$ y radio $\overline{CD}$'
$ws.Cells.Item(7, 2).WrapText = $true
$ws.Cells.Item(7, 3).Value = 'Postulado 3'
$ws.Cells.Item(7, 4).Value = 'circulo'
$ws.Cells.Item(7, 5).Value = '{''centro'':(0,1), ''radio_punto'': (1.5,-1)}'

# Row 8
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = 'Tomar un punto donde $E$ donde se interseca el circulo con la recta'
$ws.Cells.Item(8, 2).WrapText = $true
$ws.Cells.Item(8, 3).Value = '-'
$ws.Cells.Item(8, 4).Value = 'punto'
$ws.Cells.Item(8, 5).Value = '{''x'':-2.29,''y'':0,''nombre'':''E''}'
$ws.Rows.Item(8).RowHeight = 28.8

# Row 9
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = 'Tomar un punto donde $F$ donde se interseca el circulo con la recta'
$ws.Cells.Item(9, 2).WrapText = $true
$ws.Cells.Item(9, 3).Value = '-'
$ws.Cells.Item(9, 4).Value = 'punto'
$ws.Cells.Item(9, 5).Value = '{''x'':2.29,''y'':0,''nombre'':''F''}'
$ws.Rows.Item(9).RowHeight = 28.8

# Row 10
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = 'Trazar segmento $\overline{CE}$'
$ws.Cells.Item(10, 2).WrapText = $true
$ws.Cells.Item(10, 3).Value = 'Postulado 1'
$ws.Cells.Item(10, 4).Value = 'segmento'
$ws.Cells.Item(10, 5).Value = '{''x1'':-2.29, ''y1'':0, ''x2'':0, ''y2'':1}'

# Row 11
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = 'Trazar segmento $\overline{CF}$'
$ws.Cells.Item(11, 2).WrapText = $true
$ws.Cells.Item(11, 3).Value = 'Postulado 1'
$ws.Cells.Item(11, 4).Value = 'segmento'
$ws.Cells.Item(11, 5).Value = '{''x1'':2.29, ''y1'':0, ''x2'':0, ''y2'':1}'

# Row 12
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = 'Trazar segmento $\overline{EF}$'
$ws.Cells.Item(12, 2).WrapText = $true
$ws.Cells.Item(12, 3).Value = 'Postulado 1'
$ws.Cells.Item(12, 4).Value = 'segmento'
$ws.Cells.Item(12, 5).Value = '{''x1'':-2.29, ''y1'':0, ''x2'':2.29, ''y2'':0}'

# Row 13
$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = 'Bisecar $\overline{EF}$ y nombrar $G$'
$ws.Cells.Item(13, 2).WrapText = $true
$ws.Cells.Item(13, 3).Value = 'Proposición I.10'
$ws.Cells.Item(13, 4).Value = 'punto'
$ws.Cells.Item(13, 5).Value = '{''x'':0,''y'':0,''nombre'':''G''}'

# Row 14
$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = 'Trazar segmento $\overline{CG}$'
$ws.Cells.Item(14, 2).WrapText = $true
$ws.Cells.Item(14, 3).Value = '-'
$ws.Cells.Item(14, 4).Value = 'segmento'
$ws.Cells.Item(14, 5).Value = '{''x1'':0, ''y1'':1, ''x2'':0, ''y2'':0}'

# Row 15
$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 2).Value = 'Trazar el segmento $\overline{EG}$'
$ws.Cells.Item(15, 2).WrapText = $true
$ws.Cells.Item(15, 3).Value = 'Postulado 1'
$ws.Cells.Item(15, 4).Value = 'segmento'
$ws.Cells.Item(15, 5).Value = '{''x1'':-2.29, ''y1'':0, ''x2'':0, ''y2'':0}'

# Row 16
$ws.Cells.Item(16, 1).Value = 15
$ws.Cells.Item(16, 2).Value = 'Trazar el segmento $\overline{GF}$'
$ws.Cells.Item(16, 2).WrapText = $true
$ws.Cells.Item(16, 3).Value = 'Postulado 1'
$ws.Cells.Item(16, 4).Value = 'segmento'
$ws.Cells.Item(16, 5).Value = '{''x1'':2.29, ''y1'':0, ''x2'':0, ''y2'':0}'

# Row 17
$ws.Cells.Item(17, 1).Value = 16
$ws.Cells.Item(17, 2).Value = '$\overline{EG} = \overline{GF}$'
$ws.Cells.Item(17, 2).WrapText = $true
$ws.Cells.Item(17, 3).Value = 'Proposición I.10'

# Row 18
$ws.Cells.Item(18, 1).Value = 17
$ws.Cells.Item(18, 2).Value = '$\overline{CE} = \overline{CF}$'
$ws.Cells.Item(18, 2).WrapText = $true
$ws.Cells.Item(18, 3).Value = 'Definición 15'

# Row 19
$ws.Cells.Item(19, 1).Value = 18
$ws.Cells.Item(19, 2).Value = 'Por los pasos (16) y (17) decimos que $\triangle EGC = \triangle FGC$ '
$ws.Cells.Item(19, 2).WrapText = $true
$ws.Cells.Item(19, 3).Value = 'Proposición I.8'
$ws.Rows.Item(19).RowHeight = 28.8

# Row 20
$ws.Cells.Item(20, 1).Value = 19
$ws.Cells.Item(20, 2).Value = 'Como $\overline{GC}$ esta sobre $\overline{AB}$ y sus ángulos adyacentes son iguales $\overline{GC}$ es perpendicular a $\overline{AB}$'
$ws.Cells.Item(20, 2).WrapText = $true
$ws.Cells.Item(20, 3).Value = 'Definición 10'
$ws.Cells.Item(20, 4).Value = 'angulo'
$ws.Cells.Item(20, 5).Value = '{''A'':(2.29,0), ''B'':(0,0), ''C'':(0,1)}'
$ws.Rows.Item(20).RowHeight = 43.2

# Row 21
$ws.Cells.Item(21, 1).Value = 20
$ws.Cells.Item(21, 2).Value = 'Como $\overline{GC}$ esta sobre $\overline{AB}$ y sus ángulos adyacentes son iguales $\overline{GC}$ es perpendicular a $\overline{AB}$'
$ws.Cells.Item(21, 2).WrapText = $true
$ws.Cells.Item(21, 3).Value = 'Definición 10'
$ws.Cells.Item(21, 4).Value = 'angulo'
$ws.Cells.Item(21, 5).Value = '{''A'':(0,1), ''B'':(0,0), ''C'':(-2.29,0)}'
$ws.Rows.Item(21).RowHeight = 43.2

# --- Column widths ---
$ws.Columns.Item(2).ColumnWidth = 54.0
$ws.Columns.Item(3).ColumnWidth = 13.6

# --- Selection ---
$ws.Range("B5").Select()
